$d = $word.ActiveDocument

# 1. Update the iteration-review Q&A duration from 15 minutes to 10 minutes
#    ("(15' questions included)" -> "(10' questions included)").
$d.Content.Find.Execute("15", $false, $false, $false, $false, $false,
                         $true, 1, $false, "10", 2) | Out-Null

# 2. The text replace above coalesces the paragraph's runs; re-split the edited
#    run back into "(1" | "0" | "' questions included)" using zero-width,
#    bookmark-forced run boundaries (mirrors how Word splits runs around an
#    in-place keystroke edit without merging unrelated neighboring runs).
$rng2 = $d.Content.Duplicate
$rng2.Find.Execute("(10", $false, $false, $false, $false, $false,
                    $true, 1, $false, $null, 0) | Out-Null
$afterOpenParen1 = $rng2.Start + 2
$afterZero = $rng2.Start + 3
$d.Range($afterOpenParen1, $afterOpenParen1).Bookmarks.Add("TempSplitA") | Out-Null
$d.Range($afterZero, $afterZero).Bookmarks.Add("TempSplitB") | Out-Null
$d.Bookmarks.Item("TempSplitA").Delete()
$d.Bookmarks.Item("TempSplitB").Delete()

# 3. Likewise restore the trailing ", ...iteration" / "." run boundary that
#    the replace coalesced.
$rng3 = $d.Content.Duplicate
$rng3.Find.Execute("their iteration.", $false, $false, $false, $false, $false,
                    $true, 1, $false, $null, 0) | Out-Null
$beforeDot = $rng3.End - 1
$d.Range($beforeDot, $beforeDot).Bookmarks.Add("TempSplitC") | Out-Null
$d.Bookmarks.Item("TempSplitC").Delete()

# 4. Re-anchor Word's auto-tracked "last edit" (_GoBack) bookmark to the point
#    in "Systems Engineering tasks" where the author's cursor ended up
#    ("- Ensure th|at commits are done properly...").
$rng4 = $d.Content.Duplicate
$rng4.Find.Execute("- Ensure that commits are done properly", $false, $false, $false, $false, $false,
                    $true, 1, $false, $null, 0) | Out-Null
$pos = $rng4.Start + ("- Ensure th").Length
$d.Range($pos, $pos).Bookmarks.Add("_GoBack") | Out-Null
